# Updated cryptos list on Thu Apr 13 12:31:03 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on the crypto tracker sheet.
# Price/Volume cells are stored as text in this sheet, so numeric-looking
# updates are entered with a leading apostrophe (forces text) and then the
# Normal style is reapplied so no stray number-format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.221.05"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.994.89"
$ws.Range("E3").Value = "  +6.22%  "

$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").Value = "'324.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").Value = "'0.5086"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("D8").Value = "'0.4115"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.95%  "

$ws.Range("D9").Value = "'0.08707"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.09%  "

$ws.Range("D10").Value = "'1.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.54%  "

$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").Value = "'24.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.53%  "

$ws.Range("D13").Value = "1.986.69"
$ws.Range("E13").Value = "  +6.01%  "

$ws.Range("D14").Value = "'6.505"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("D15").Value = "'7.398"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.94%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "'93.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.57%  "

$ws.Range("D18").Value = "'0.00001115"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("E19").Value = "  +1.62%  "

$ws.Range("D20").Value = "'18.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").Value = "'6.092"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.33%  "

$ws.Range("D23").Value = "30.283.75"
$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").Value = "'11.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.07%  "

$ws.Range("D25").Value = "'2.210"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").Value = "2.223.11"
$ws.Range("E26").Value = "  +6.35%  "

$ws.Range("D27").Value = "'22.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.34%  "

$ws.Range("D28").Value = "'163.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "

$ws.Range("D29").Value = "'2.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.13%  "

$ws.Range("D30").Value = "'130.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.76%  "

$ws.Range("D31").Value = "'1.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.36%  "

$ws.Range("D32").Value = "'0.1051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.68%  "

$ws.Range("D33").Value = "'6.058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").Value = "'3.811"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "

$ws.Range("D35").Value = "'1.313"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.18%  "

$ws.Range("D36").Value = "'0.02478"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "

$ws.Range("D37").Value = "'5.385"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.08%  "

$ws.Range("D38").Value = "'0.06535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("D39").Value = "'0.2191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "

$ws.Range("D40").Value = "'8.858"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.32%  "

$ws.Range("D41").Value = "'0.6579"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.55%  "

$ws.Range("D42").Value = "'11.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.52%  "

$ws.Range("D43").Value = "'1.222"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").Value = "'13.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.59%  "

$ws.Range("D45").Value = "'0.6116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.81%  "

$ws.Range("D46").Value = "'2.194"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.80%  "

$ws.Range("D47").Value = "'3.658"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").Value = "'123.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").Value = "'1.224"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("D50").Value = "'79.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.73%  "

$ws.Range("D51").Value = "'0.06878"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "

